# Scheduled-runner refresh of market/profit figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
# Values below are the newly refreshed figures; a couple of rows flip a
# NQ/HQ average price to/from zero, which drops (ClearContents) or creates
# the corresponding LeveProfit cell, matching Excel's own behavior when the
# source data crosses a zero boundary.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 385.22223
$ws.Range("I12").Value = 316.7143
$ws.Range("J12").Value = 625
$ws.Range("K12").Value = 316.7143
$ws.Range("L12").Value = 625
$ws.Range("M12").Value = -146.7143
$ws.Range("N12").Value = -965

$ws.Range("H19").Value = 1977.8148
$ws.Range("I19").Value = 1810.2354
$ws.Range("J19").Value = 2262.7
$ws.Range("K19").Value = 1810.2354
$ws.Range("L19").Value = 2262.7
$ws.Range("M19").Value = -1635.2354
$ws.Range("N19").Value = -2612.7

$ws.Range("H98").Value = 4469.067
$ws.Range("I98").Value = 5245.9
$ws.Range("K98").Value = 5245.9
$ws.Range("M98").Value = -3747.9

$ws.Range("H116").Value = 4261.4707
$ws.Range("I116").Value = 4386.273
$ws.Range("J116").Value = 4032.6667
$ws.Range("K116").Value = 4386.273
$ws.Range("L116").Value = 4032.6667
$ws.Range("M116").Value = -944.2730000000001
$ws.Range("N116").Value = -10916.6667

$ws.Range("H122").Value = 4469.067
$ws.Range("I122").Value = 5245.9
$ws.Range("K122").Value = 15737.7
$ws.Range("M122").Value = -13287.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2486.3428
$ws.Range("I2").Value = 737.3200000000001
$ws.Range("K2").Value = 737.3200000000001
$ws.Range("M2").Value = -624.3200000000001

$ws.Range("H32").Value = 8542.177
$ws.Range("I32").Value = 6414.7334
$ws.Range("J32").Value = 24498
$ws.Range("K32").Value = 6414.7334
$ws.Range("L32").Value = 24498
$ws.Range("M32").Value = -6127.7334
$ws.Range("N32").Value = -25072

$ws.Range("H45").Value = 2216
$ws.Range("I45").Value = 2418.7778
$ws.Range("J45").Value = 1911.8334
$ws.Range("K45").Value = 2418.7778
$ws.Range("L45").Value = 1911.8334
$ws.Range("M45").Value = -2041.7778
$ws.Range("N45").Value = -2665.8334

$ws.Range("H74").Value = 15481.104
$ws.Range("J74").Value = 29044.133
$ws.Range("L74").Value = 29044.133
$ws.Range("N74").Value = -30792.133

$ws.Range("H77").Value = 15481.104
$ws.Range("J77").Value = 29044.133
$ws.Range("L77").Value = 145220.665
$ws.Range("N77").Value = -153956.665

$ws.Range("H112").Value = 28386.572
$ws.Range("J112").Value = 28386.572
$ws.Range("L112").Value = 28386.572
$ws.Range("N112").Value = -31340.572

$ws.Range("H116").Value = 2486.3428
$ws.Range("I116").Value = 737.3200000000001
$ws.Range("K116").Value = 737.3200000000001
$ws.Range("M116").Value = 1556.68

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2486.3428
$ws.Range("I3").Value = 737.3200000000001
$ws.Range("K3").Value = 737.3200000000001
$ws.Range("M3").Value = -623.3200000000001

$ws.Range("H9").Value = 4014800
$ws.Range("J9").Value = 4014800
$ws.Range("L9").Value = 4014800
$ws.Range("N9").Value = -4015136

$ws.Range("H22").Value = 1175
$ws.Range("I22").Value = 1175
$ws.Range("K22").Value = 1175
$ws.Range("M22").Value = -1002

$ws.Range("H94").Value = 2823.8147
$ws.Range("I94").Value = 2708.558
$ws.Range("J94").Value = 3274.3635
$ws.Range("K94").Value = 2708.558
$ws.Range("L94").Value = 3274.3635
$ws.Range("M94").Value = -2257.558
$ws.Range("N94").Value = -4176.363499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 437142.84
$ws.Range("I4").Value = 15000.75
$ws.Range("J4").Value = 999999
$ws.Range("K4").Value = 15000.75
$ws.Range("L4").Value = 999999
$ws.Range("M4").Value = -14888.75
$ws.Range("N4").Value = -1000223

$ws.Range("H13").Value = 4000
$ws.Range("J13").Value = 4000
$ws.Range("L13").Value = 4000
$ws.Range("N13").Value = -4278

$ws.Range("H31").Value = 26266.809
$ws.Range("I31").Value = 12216.4
$ws.Range("K31").Value = 12216.4
$ws.Range("M31").Value = -11921.4

$ws.Range("H34").Value = 26266.809
$ws.Range("I34").Value = 12216.4
$ws.Range("K34").Value = 12216.4
$ws.Range("M34").Value = -12014.4

$ws.Range("H86").Value = 10639.077
$ws.Range("I86").Value = 14479
$ws.Range("J86").Value = 7347.7144
$ws.Range("K86").Value = 14479
$ws.Range("L86").Value = 7347.7144
$ws.Range("M86").Value = -13356
$ws.Range("N86").Value = -9593.714400000001

$ws.Range("H89").Value = 10639.077
$ws.Range("I89").Value = 14479
$ws.Range("J89").Value = 7347.7144
$ws.Range("K89").Value = 72395
$ws.Range("L89").Value = 36738.572
$ws.Range("M89").Value = -66779
$ws.Range("N89").Value = -47970.572

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 79010530
$ws.Range("I4").Value = 102613660
$ws.Range("J4").Value = 333397.34
$ws.Range("K4").Value = 307840980
$ws.Range("L4").Value = 1000192.02
$ws.Range("M4").Value = -307840868
$ws.Range("N4").Value = -1000416.02

$ws.Range("H14").Value = 254.75
$ws.Range("I14").Value = 254.75
$ws.Range("K14").Value = 764.25
$ws.Range("M14").Value = -591.25

$ws.Range("H80").Value = 14637.846
$ws.Range("J80").Value = 21649.666
$ws.Range("L80").Value = 64948.99800000001
$ws.Range("N80").Value = -66820.99800000001

$ws.Range("H83").Value = 14637.846
$ws.Range("J83").Value = 21649.666
$ws.Range("L83").Value = 194846.994
$ws.Range("N83").Value = -204206.994

$ws.Range("H114").Value = 844.1667
$ws.Range("I114").Value = 913
$ws.Range("J114").Value = 500
$ws.Range("K114").Value = 2739
$ws.Range("L114").Value = 1500
$ws.Range("M114").Value = 515
$ws.Range("N114").Value = -8008

$ws.Range("H134").Value = 6555.1562
$ws.Range("I134").Value = 2263.3333
$ws.Range("K134").Value = 6789.999899999999
$ws.Range("M134").Value = -1719.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 19990
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 19990
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 19990
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -21630

$ws.Range("H113").Value = 3723.1538
$ws.Range("I113").Value = 1987.625
$ws.Range("K113").Value = 1987.625
$ws.Range("M113").Value = 182.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 20000000
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()

$ws.Range("H7").Value = 13126.096
$ws.Range("I7").Value = 14961
$ws.Range("K7").Value = 14961
$ws.Range("M7").Value = -14849

$ws.Range("H22").Value = 5677.45
$ws.Range("I22").Value = 2235.6365
$ws.Range("K22").Value = 2235.6365
$ws.Range("M22").Value = -1940.6365

$ws.Range("H27").Value = 5677.45
$ws.Range("I27").Value = 2235.6365
$ws.Range("K27").Value = 2235.6365
$ws.Range("M27").Value = -2128.6365

$ws.Range("H110").Value = 46276
$ws.Range("J110").Value = 46276
$ws.Range("L110").Value = 46276
$ws.Range("N110").Value = -54456

$ws.Range("H126").Value = 13126.096
$ws.Range("I126").Value = 14961
$ws.Range("K126").Value = 44883
$ws.Range("M126").Value = -42413

$ws.Range("H132").Value = 959655.5
$ws.Range("I132").Value = 2045.9697
$ws.Range("K132").Value = 6137.909100000001
$ws.Range("M132").Value = -3607.909100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 1905.7059
$ws.Range("I2").Value = 2294.7856
$ws.Range("K2").Value = 2294.7856
$ws.Range("M2").Value = -2182.7856

$ws.Range("H26").Value = 9499.5
$ws.Range("I26").Value = 9499.5
$ws.Range("K26").Value = 9499.5
$ws.Range("M26").Value = -9206.5

$ws.Range("H126").Value = 14546.944
$ws.Range("I126").Value = 8856
$ws.Range("J126").Value = 20237.889
$ws.Range("K126").Value = 26568
$ws.Range("L126").Value = 60713.667
$ws.Range("M126").Value = -24098
$ws.Range("N126").Value = -65653.667
